$d = $word.ActiveDocument

# --- Edit 1: "polynomial runtime" -> "runtime" ---
$null = $d.Content.Find.Execute(
    "We gain insight to the polynomial runtime of the algorithm",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We gain insight to the runtime of the algorithm", 2)

# --- Edit 2: rewrite the "Conclusions and Future Work" closing paragraph ---
# Locate the paragraph that currently starts the conclusion text.
$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "As we adventure into the future*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "could not locate target paragraph"
}

$fragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="2617405A" w14:textId="1FC98346" w:rsidR="001104EA" w:rsidRPr="001104EA" w:rsidRDefault="001104EA" w:rsidP="00E81C91" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>The distributed nature of the internet as well a</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>distributed algorithms</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> for the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>steiner</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> tree problem, which </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">are inherently scalable and </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>could be a promising avenue for future work.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> Also, because the general problem is NP </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">hard (by reduction to 3SAT) </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">approximation algorithms under problem specific assumptions are and will continue to be </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>an</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>fruitful</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> avenue </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>for</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> research. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>As we adventure into the future, we stand to see the expansion of the internet</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> to what </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>is project</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> to be as much as 5 billion people as soon as 2020 (by some </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>predications</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> as well as the rise of the </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">internet of things, and multi-robotic networked systems. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">Network demands will become heavier, and more dynamic than ever. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">Our civilization will see some of it’s the most wonderful implications of these technologies as well as the great challenges that come with them in the coming decades. </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve">This trajectory of growth implies that </w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>scalable</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> and</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> flexible</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> dynamic multipoint communication</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> is and will increasingly become a fundamentally important theoretical and practical problem in computer science</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr></w:p>
'@

$target = $d.Paragraphs($targetIndex).Range
$target.InsertXML($fragment)
